# IMP: Add Tax No. In Report Gl Detail
#
# The report template's line-item table gets a new "Tax No." column
# inserted right after the "Number" column (between "Number" and
# "Description"), both in the header row (11) and in the per-line
# template row (15), which holds the Handlebars placeholder
# "{{tax_no}}".
#
# Inserting a whole column at D shifts every existing column from D
# onward one position to the right (D->E, E->F, F..I->G..J, J->K,
# K->L), extends the merged header cells (C2:E2 etc.) to C:F, and
# carries formatting the same way Excel's own "Insert Column" does.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D - everything from D rightwards shifts
# right by one, merged cells grow and formatting is carried across,
# matching the dimension change from K18 to L18.
$ws.Columns("D:D").Insert()

# New column header + template placeholder.
$ws.Range("D11").Value = "Tax No."
$ws.Range("D15").Value = "{{tax_no}}"

# Give the new column its own width (close to the ~15.69-character
# width used in the final report; Excel quantizes ColumnWidth to
# whole pixels so this is the nearest representable value).
$ws.Columns("D:D").ColumnWidth = 14.8333

# Restore the view to show the sheet from the top-left and leave the
# selection on the new column's template cell.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E15").Select()
